$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Activate()

# --- Add the new test case row (row 20) -----------------------------------
# Write cells in the same order the original author's shared-string table
# grew in (Description, Jira id, TCID, Runmode, Results) so new shared
# strings land at the same indices as the target file.
$ws.Range("C20").Value = "Verify that follower of the article is able to start conversation from home page when some one commented on the article he is following."
$ws.Range("B20").Value = "OPQA-1012"
$ws.Range("A20").Value = "TestCase_F19"
$ws.Range("D20").Value = "Y"
$ws.Range("E20").Value = "PASS"

# Row 2's Results cell moves from SKIP to PASS now that the feature works.
$ws.Range("E2").Value = "PASS"

# --- Formatting for the new row --------------------------------------------
# Row 20 should look like row 19 (borders + shaded fill), so copy formats.
$ws.Range("A19:E19").Copy() | Out-Null
$ws.Range("A20:E20").PasteSpecial(-4122) | Out-Null
# ...except column B, which should match the plain bordered style used in
# column C (no shading), same as the rest of the "odd one out" rows.
$ws.Range("C19").Copy() | Out-Null
$ws.Range("B20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- View state --------------------------------------------------------
$excel.Goto($ws.Range("D13"), $true)
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 2
$win.TopLeftCell = $ws.Range("B1")
$ws.Range("D13").Select()
